$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colIndices = @(2, 3, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.022865631252103, 1.025510287619976, 1.032608628337462, 1.041416197227252, 1, 1.026706965650817, 1.028049065184334, 1.028335647395785, 1.0354133672848, 1.044195774322585, 1.013324148161443),
    @(1.02, 1.023998909902228, 1.026556119376375, 1.0336159513703, 1.042513387903579, 1, 1.026675851733607, 1.028819623087733, 1.029188150035785, 1.036228967212125, 1.045102835912574, 1.013586760982975),
    @(1.02, 1.024732554952344, 1.027233488926321, 1.034268178707507, 1.043223275939709, 1, 1.026653328498527, 1.029318058496989, 1.029739827277754, 1.036756530907662, 1.045689109917733, 1.013756436762798),
    @(1.02, 1.025041060988944, 1.027518409937326, 1.034542475959977, 1.043521696316437, 1, 1.026643286458578, 1.029527560827157, 1.029971764690024, 1.036978274530145, 1.045935422132413, 1.013827707973501),
    @(1.02, 1.025092865320127, 1.027566258511916, 1.034588537587128, 1.043571801429935, 1, 1.026641566709819, 1.029562734851993, 1.030010708725819, 1.037015503667686, 1.045976769795241, 1.013839671173715),
    @(1.02, 1.024736676903317, 1.027237295446253, 1.034271843484629, 1.043227263512913, 1, 1.026653196570139, 1.029320858034129, 1.029742926388387, 1.036759494030453, 1.045692401772597, 1.013757389329996),
    @(1.02, 1.023248558542754, 1.025863597554226, 1.032948970414474, 1.041787011763677, 1, 1.026696944921733, 1.028309513480885, 1.028623744138844, 1.03568904113286, 1.044502455489732, 1.013412951574591),
    @(1.02, 1.020628858867363, 1.023447913184159, 1.030621136565004, 1.039248600182973, 1, 1.026755770427594, 1.026526111331389, 1.026651988622907, 1.033801361307025, 1.042400616899688, 1.01280407996468),
    @(1.02, 1.018884060854426, 1.021840773836129, 1.029071428695964, 1.037556003022065, 1, 1.026782765030883, 1.025336305053232, 1.025337736353184, 1.032541966270522, 1.040996051377452, 1.012396873390566),
    @(1.02, 1.018128927993625, 1.02114564853345, 1.028400905471588, 1.0368230154695, 1, 1.026791566969526, 1.024820896454931, 1.024768707518954, 1.031996410648891, 1.040387070648639, 1.012220242574225),
    @(1.02, 1.017848493432861, 1.020887564372235, 1.028151919877941, 1.036550739015063, 1, 1.026794403602652, 1.024629418188923, 1.024557352217143, 1.031793732427587, 1.040160748898479, 1.012154587810452),
    @(1.02, 1.017908645108707, 1.02094291904258, 1.02820532466907, 1.036609143791983, 1, 1.026793814705538, 1.024670492405683, 1.024602688319847, 1.031837209160635, 1.040209301072949, 1.012168673071258),
    @(1.02, 1.018105746061512, 1.0211243128452, 1.028380322683049, 1.036800509262284, 1, 1.026791810270159, 1.024805069467873, 1.024751236675883, 1.031979657909915, 1.040368365252613, 1.01221481647366),
    @(1.02, 1.018227193772282, 1.021236090957537, 1.028488154880268, 1.036918414223737, 1, 1.026790517948808, 1.024887982512903, 1.024842763183932, 1.032067420700659, 1.040466354154753, 1.012243240823324),
    @(1.02, 1.018934184366349, 1.021886923345777, 1.029115939860446, 1.037604647272606, 1, 1.026782120130775, 1.025370506475296, 1.025375501970749, 1.032578168159192, 1.041036450697671, 1.012408589314631),
    @(1.02, 1.019377760819289, 1.022295380986882, 1.029509869392943, 1.038035081268085, 1, 1.026776080152916, 1.025673123322304, 1.02570968837795, 1.032898485076752, 1.041393844596261, 1.012512225688069),
    @(1.02, 1.019636527831213, 1.022533702376751, 1.029739691088815, 1.038286138238406, 1, 1.026772278530361, 1.025849613922155, 1.025904618605052, 1.033085298403886, 1.041602229732163, 1.012572645343663),
    @(1.02, 1.019724766898286, 1.022614976536948, 1.029818062707429, 1.038371740836364, 1, 1.026770934990617, 1.025909789156569, 1.025971085624639, 1.033148993174146, 1.041673270666647, 1.012593241848927),
    @(1.02, 1.019330165509604, 1.022251549596712, 1.029467599387232, 1.037988900590215, 1, 1.026776756996141, 1.025640657540494, 1.025673832829431, 1.03286412039084, 1.041355507542992, 1.012501109556609),
    @(1.02, 1.01804770322297, 1.021070893683116, 1.028328788016373, 1.036744157238475, 1, 1.02679241246641, 1.02476544077879, 1.024707492687441, 1.03193771125911, 1.040321528105108, 1.012201229667392),
    @(1.02, 1.017241688295548, 1.020329241304278, 1.027613215159834, 1.035961466887169, 1, 1.026799752091084, 1.024214968699833, 1.024099958608288, 1.031355040977366, 1.03967073466997, 1.012012415903474),
    @(1.02, 1.017668942096462, 1.020722341602748, 1.027992511911882, 1.036376392547362, 1, 1.026796098222743, 1.024506802344464, 1.024422020017272, 1.031663944528271, 1.040015797881986, 1.012112534992648),
    @(1.02, 1.019351671668581, 1.022271354883401, 1.029486699230912, 1.038009767674897, 1, 1.026776452020753, 1.025655327493544, 1.025690034392924, 1.032879648381331, 1.0413728306439, 1.012506132549609),
    @(1.02, 1.021305816340776, 1.024071840202923, 1.031222553261458, 1.039904898468015, 1, 1.026742721087137, 1.026987316529266, 1.027161689745313, 1.034289537990128, 1.04294458285845, 1.012961715934136)
)

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $colIndices[$c]).Value = $rowVals[$c]
    }
}
